# Bug Report and Test case on E-commerce site
# Renumber Bug IDs from "DARAZ-BUG-0NN" to "BUG-0NN"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    $new = $old -replace "^DARAZ-", ""
    $cell.Value = $new
}

# Widen column A (Bug ID) and column G (Expected Result)
$ws.Columns.Item(1).ColumnWidth = 21.166666666666668
$ws.Columns.Item(7).ColumnWidth = 43

# Adjust zoom and select the next empty row below the data
$win = $wb.Windows.Item(1)
$win.Zoom = 115
[void]$ws.Range("A17").Select()

# Resize the application window
$win.Width = 20490
$win.Height = 7320
$win.Left = 0
$win.Top = 0
